# Rename the "_old" / "_new" header suffixes to the respective format-version
# suffixes ("_FV2210" for the old/left-hand side, "_FV2304" for the new/
# right-hand side), matching the commit's move from generic "old"/"new"
# column names to names carrying the concrete <formatversion> they refer to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Segmentname_FV2210"
$ws.Range("B1").Value = "Segmentgruppe_FV2210"
$ws.Range("C1").Value = "Segment_FV2210"
$ws.Range("D1").Value = "Datenelement_FV2210"
$ws.Range("E1").Value = "Segment ID_FV2210"
$ws.Range("F1").Value = "Code_FV2210"
$ws.Range("G1").Value = "Qualifier_FV2210"
$ws.Range("H1").Value = "Beschreibung_FV2210"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2210"
$ws.Range("J1").Value = "Bedingung_FV2210"

# K1 ("diff") is untouched.

$ws.Range("L1").Value = "Segmentname_FV2304"
$ws.Range("M1").Value = "Segmentgruppe_FV2304"
$ws.Range("N1").Value = "Segment_FV2304"
$ws.Range("O1").Value = "Datenelement_FV2304"
$ws.Range("P1").Value = "Segment ID_FV2304"
$ws.Range("Q1").Value = "Code_FV2304"
$ws.Range("R1").Value = "Qualifier_FV2304"
$ws.Range("S1").Value = "Beschreibung_FV2304"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("U1").Value = "Bedingung_FV2304"

# Turn the used range into a real Excel Table ("Table1") covering the
# whole sheet, with column headers taken from row 1 (which we just
# relabeled above).
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U82"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (split/freeze after row 1, top-left cell of the
# scrolling pane is A2).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
